$wb = $excel.ActiveWorkbook

# --- "Test Server" sheet: add Tes12 / t24tes12 row ---
$wsTestServer = $wb.Worksheets.Item("Test Server")
$wsTestServer.Range("A6").Value = "Tes12"
$wsTestServer.Range("B6").Value = "t24tes12"
$wsTestServer.Range("A7").Select()

# --- "Server IP" sheet: add Tes12 / 10.169.1.36 row ---
$wsServerIP = $wb.Worksheets.Item("Server IP")
$wsServerIP.Range("A6").Value = "Tes12"
$wsServerIP.Range("B6").Value = "10.169.1.36"
$wsServerIP.Range("A1").Select()

# --- "Users" sheet: no longer the active/selected tab ---
$wsUsers = $wb.Worksheets.Item("Users")

# --- New "Fechas" sheet with date parametrization data ---
$wsFechas = $wb.Worksheets.Add([System.Type]::Missing, $wb.Worksheets.Item($wb.Worksheets.Count))
$wsFechas.Name = "Fechas"

$wsFechas.Range("A2").Value = "COB"
$wsFechas.Range("A3").Value = "DIA"
$wsFechas.Range("A4").Value = "VALOR"

$wsFechas.Range("D1").Value = "TES10"
$wsFechas.Range("E1").Value = "TES11"
$wsFechas.Range("F1").Value = "TES12"

$wsFechas.Range("B1").Value = 702
$wsFechas.Range("C1").Value = 708

$wsFechas.Range("B2:F4").NumberFormat = "@"
$wsFechas.Range("B2:F4").HorizontalAlignment = -4152

$wsFechas.Range("B2").Value = "20230612"
$wsFechas.Range("C2").Value = "20230901"
$wsFechas.Range("D2").Value = "20230901"
$wsFechas.Range("E2").Value = "20230829"

$wsFechas.Range("B3").Value = "20230612"
$wsFechas.Range("C3").Value = "20230901"
$wsFechas.Range("D3").Value = "20230901"
$wsFechas.Range("E3").Value = "20230829"

$wsFechas.Range("B4").Value = "20230612"
$wsFechas.Range("C4").Value = "20230901"
$wsFechas.Range("D4").Value = "20230901"
$wsFechas.Range("E4").Value = "20230829"

$wsFechas.PageSetup.Orientation = 1

$wsFechas.Range("F2").Select()
$wsFechas.Activate()
